$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed numeric values (rows 2-6) ---
$ws.Range("D2").Value = 340
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 27
$ws.Range("G2").Value = -4
$ws.Range("H2").Value = -10
$ws.Range("I2").Value = -10
$ws.Range("K2").Value = 617
$ws.Range("L2").Value = 43
$ws.Range("M2").Value = 575
$ws.Range("N2").Value = 575
$ws.Range("P2").Value = 79
$ws.Range("Q2").Value = 65
$ws.Range("R2").Value = 11
$ws.Range("S2").Value = -5
$ws.Range("T2").Value = 6
$ws.Range("U2").Value = 60
$ws.Range("W2").Value = 8
$ws.Range("X2").Value = -2.92
$ws.Range("Y2").Value = -1.72
$ws.Range("Z2").Value = -1.56
$ws.Range("AA2").Value = 7.41
$ws.Range("AB2").Value = 572.99
$ws.Range("AC2").Value = -633
$ws.Range("AD2").Value = -30.32
$ws.Range("AE2").Value = 41901
$ws.Range("AF2").Value = 0.46
$ws.Range("AG2").Value = 350
$ws.Range("AH2").Value = 1.82
$ws.Range("AI2").Value = -48.26
$ws.Range("AJ2").Value = 1570797
$ws.Range("D3").Value = 283
$ws.Range("E3").Value = 17
$ws.Range("F3").Value = 17
$ws.Range("G3").Value = 21
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 13
$ws.Range("K3").Value = 642
$ws.Range("L3").Value = 46
$ws.Range("M3").Value = 596
$ws.Range("N3").Value = 596
$ws.Range("P3").Value = 79
$ws.Range("Q3").Value = 22
$ws.Range("R3").Value = 9
$ws.Range("S3").Value = -5
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 21
$ws.Range("W3").Value = 5.9
$ws.Range("X3").Value = 4.71
$ws.Range("Y3").Value = 2.28
$ws.Range("Z3").Value = 2.12
$ws.Range("AA3").Value = 7.74
$ws.Range("AB3").Value = 583.54
$ws.Range("AC3").Value = 850
$ws.Range("AD3").Value = 23.11
$ws.Range("AE3").Value = 43467
$ws.Range("AF3").Value = 0.45
$ws.Range("AG3").Value = 350
$ws.Range("AH3").Value = 1.78
$ws.Range("AI3").Value = 35.95
$ws.Range("AJ3").Value = 1570797
$ws.Range("D4").Value = 266
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 44
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 30
$ws.Range("K4").Value = 643
$ws.Range("L4").Value = 46
$ws.Range("M4").Value = 597
$ws.Range("N4").Value = 597
$ws.Range("P4").Value = 79
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = 33
$ws.Range("S4").Value = -5
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 12
$ws.Range("W4").Value = 2.56
$ws.Range("X4").Value = 11.11
$ws.Range("Y4").Value = 4.95
$ws.Range("Z4").Value = 4.59
$ws.Range("AA4").Value = 7.62
$ws.Range("AB4").Value = 675.35
$ws.Range("AC4").Value = 1880
$ws.Range("AD4").Value = 11.99
$ws.Range("AE4").Value = 43555
$ws.Range("AF4").Value = 0.52
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 1.11
$ws.Range("AI4").Value = 11.61
$ws.Range("AJ4").Value = 1570797
$ws.Range("D5").Value = 220
$ws.Range("E5").Value = -7
$ws.Range("F5").Value = -7
$ws.Range("G5").Value = -30
$ws.Range("H5").Value = -30
$ws.Range("I5").Value = -30
$ws.Range("K5").Value = 599
$ws.Range("L5").Value = 38
$ws.Range("M5").Value = 561
$ws.Range("N5").Value = 561
$ws.Range("P5").Value = 79
$ws.Range("Q5").Value = -24
$ws.Range("R5").Value = 15
$ws.Range("S5").Value = -3
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = -26
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = -3.27
$ws.Range("X5").Value = -13.53
$ws.Range("Y5").Value = -5.15
$ws.Range("Z5").Value = -4.8
$ws.Range("AA5").Value = 6.77
$ws.Range("AB5").Value = 629.1900000000001
$ws.Range("AC5").Value = -1897
$ws.Range("AD5").Value = -10.39
$ws.Range("AE5").Value = 40868
$ws.Range("AF5").Value = 0.48
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 1.27
$ws.Range("AI5").Value = -11.51
$ws.Range("AJ5").Value = 1570797
$ws.Range("D6").Value = 185
$ws.Range("E6").Value = -17
$ws.Range("F6").Value = -17
$ws.Range("G6").Value = -7
$ws.Range("H6").Value = -7
$ws.Range("I6").Value = -7
$ws.Range("K6").Value = 574
$ws.Range("L6").Value = 28
$ws.Range("M6").Value = 546
$ws.Range("N6").Value = 546
$ws.Range("P6").Value = 79
$ws.Range("Q6").Value = -2
$ws.Range("R6").Value = 11
$ws.Range("S6").Value = -11
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = -6
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = -9.43
$ws.Range("X6").Value = -3.69
$ws.Range("Y6").Value = -1.23
$ws.Range("Z6").Value = -1.16
$ws.Range("AA6").Value = 5.08
$ws.Range("AB6").Value = 607.75
$ws.Range("AC6").Value = -433
$ws.Range("AD6").Value = -50.19
$ws.Range("AE6").Value = 40853
$ws.Range("AF6").Value = 0.53
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 1.15
$ws.Range("AI6").Value = -49.13
$ws.Range("AJ6").Value = 1570797

# --- Clear cells removed from rows 2-5 (columns J, O, V no longer populated) ---
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Rows 7-9: all metric columns cleared, only A/B/C (index/label/period) remain ---
$ws.Range("D7:AJ9").ClearContents()
